$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.564.99"
$ws.Range("E2").Value = "  +0.52%  "

$ws.Range("D3").Value = "1.953.55"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'244.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").Value = "'0.614"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.66%  "

$ws.Range("D7").Value = "'58.38"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.81%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.376"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.15%  "

$ws.Range("D10").Value = "'0.0790"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.18%  "

$ws.Range("E11").Value = "  -0.98%  "

$ws.Range("D12").Value = "'0.841"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.62%  "

$ws.Range("D13").Value = "'14.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.84%  "

$ws.Range("D14").Value = "2.238.33"
$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("D15").Value = "'21.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.35%  "

$ws.Range("D16").Value = "'5.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.99%  "

$ws.Range("D17").Value = "1.955.70"
$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("D18").Value = "36.467.21"
$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("D19").Value = "'69.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.69%  "

$ws.Range("D20").Value = "0.0₃0849"
$ws.Range("E20").Value = "  -2.39%  "

$ws.Range("D21").Value = "'229.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.74%  "

$ws.Range("D22").Value = "'5.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("E24").Value = "  +2.98%  "

$ws.Range("D25").Value = "'2.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.30%  "

$ws.Range("D26").Value = "'9.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").Value = "'160.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.42%  "

$ws.Range("D28").Value = "'0.136"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.69%  "

$ws.Range("D29").Value = "'19.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.69%  "

$ws.Range("D30").Value = "'0.120"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.79%  "

$ws.Range("D31").Value = "'1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.40%  "

$ws.Range("D32").Value = "'4.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.24%  "

$ws.Range("D33").Value = "'0.0611"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.25%  "

$ws.Range("D34").Value = "'4.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.39%  "

$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "'3.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +15.15%  "

$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.31%  "

$ws.Range("D37").Value = "'2.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.29%  "

$ws.Range("E38").Value = "  -1.40%  "

$ws.Range("D39").Value = "'5.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -14.87%  "

$ws.Range("D40").Value = "'0.0976"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.10%  "

$ws.Range("D41").Value = "'2.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.43%  "

$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("D43").Value = "'0.0210"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").Value = "'15.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.53%  "

$ws.Range("D45").Value = "1.364.88"
$ws.Range("E45").Value = "  +1.81%  "

$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'1.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.51%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'87.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.94%  "

$ws.Range("D48").Value = "'7.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.15%  "

$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").Value = "2.129.13"
$ws.Range("E50").Value = "  +0.48%  "

$ws.Range("D51").Value = "'44.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.33%  "
